$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New test case rows: sc17, sc18, sc19 (hierarchal-condition tests added for PR#1420)
# The text cells are written in a specific order so freshly-introduced shared
# strings land in the same sequence the original authoring produced.
$ws.Cells.Item(21, 1).Value = "sc17"
$ws.Cells.Item(21, 5).Value = "Hierarchal conditions. All locations duplicated some with blank cond tag. Tests PR#1420"

$ws.Cells.Item(22, 1).Value = "sc18"
$ws.Cells.Item(22, 5).Value = "Simple two location version of sc17. Includes duplicate loc 2 with blank cond tag."

$ws.Cells.Item(23, 5).Value = "Simple two location version of sc17. No duplicate loc 2."
$ws.Cells.Item(23, 1).Value = "sc19"

# Remaining numeric / reused-text columns for the three new rows.
$ws.Cells.Item(21, 2).Value = 1
$ws.Cells.Item(21, 3).Value = 2
$ws.Cells.Item(21, 4).Value = 10
$ws.Cells.Item(21, 6).Value = "complete"
$ws.Cells.Item(21, 7).Value = "yes"
$ws.Cells.Item(21, 8).Value = "done"

$ws.Cells.Item(22, 2).Value = 1
$ws.Cells.Item(22, 3).Value = 2
$ws.Cells.Item(22, 4).Value = 2
$ws.Cells.Item(22, 6).Value = "complete"
$ws.Cells.Item(22, 7).Value = "yes"
$ws.Cells.Item(22, 8).Value = "done"

$ws.Cells.Item(23, 2).Value = 1
$ws.Cells.Item(23, 3).Value = 2
$ws.Cells.Item(23, 4).Value = 2
$ws.Cells.Item(23, 6).Value = "complete"
$ws.Cells.Item(23, 7).Value = "yes"
$ws.Cells.Item(23, 8).Value = "done"
